$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold numeric-looking text values (prices, percentages).
# Temporarily force Text format so Excel does not auto-convert these to numbers,
# then restore the original (default/"Normal") style so cell formatting is unchanged.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '51.948.80'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '2.934.57'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '353.70'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").Value = '107.87'
$ws.Range("E6").Value = '  -4.31%  '
$ws.Range("D7").Value = '0.563'
$ws.Range("E7").Value = '  +1.42%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").Value = '38.22'
$ws.Range("E10").Value = '  -3.42%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '0.0864'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("E13").Value = '  -2.90%  '
$ws.Range("D14").Value = '3.435.39'
$ws.Range("E14").Value = '  +2.18%  '
$ws.Range("D15").Value = '7.76'
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").Value = '2.996.20'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("D17").Value = '0.970'
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '51.990.35'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '3.45'
$ws.Range("E19").Value = '  +4.23%  '
$ws.Range("D20").Value = '7.52'
$ws.Range("E20").Value = '  -0.83%  '
$ws.Range("D21").Value = '13.61'
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").Value = '0.0₃0973'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '69.89'
$ws.Range("E23").Value = '  -1.54%  '
$ws.Range("D24").Value = '265.66'
$ws.Range("E24").Value = '  -1.24%  '
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").Value = '0.175'
$ws.Range("E26").Value = '  -3.58%  '
$ws.Range("D27").Value = '26.84'
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").Value = '7.60'
$ws.Range("E28").Value = '  +14.32%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").Value = '0.105'
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("E31").Value = '  -2.98%  '
$ws.Range("D32").Value = '36.57'
$ws.Range("E32").Value = '  -2.12%  '
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("E34").Value = '  -2.42%  '
$ws.Range("D35").Value = '52.09'
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").Value = '0.0435'
$ws.Range("E36").Value = '  -3.12%  '
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("E38").Value = '  -3.93%  '
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("E40").Value = '  -4.74%  '
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("D43").Value = '23.15'
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("D44").Value = '118.09'
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("D45").Value = '2.16'
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("E46").Value = '  -3.68%  '
$ws.Range("D47").Value = '2.121.48'
$ws.Range("E47").Value = '  -2.14%  '
$ws.Range("D48").Value = '3.37'
$ws.Range("E48").Value = '  -3.09%  '
$ws.Range("D49").Value = '0.0346'
$ws.Range("E49").Value = '  +1.14%  '
$ws.Range("E50").Value = '  -8.69%  '
$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").Value = '0.905'
$ws.Range("E51").Value = '  -5.03%  '

$priceRange.Style = "Normal"
